# ----------------------------------------------------------------------------
# "Refined metadata to be additional tab"
#
# The panel was re-queried; this bumps the per-gene "time_taken" timestamps
# on the "data" sheet, and moves the single-row summary metadata (panel
# name/id/version/query info) that used to live bolted onto "data" into its
# own dedicated "metadata" worksheet.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Refresh the "time_taken" (column F) timestamps on the data sheet ------
$data.Range("F2").Value = "2021-10-05 14:22:50.678376"
$data.Range("F3").Value = "2021-10-05 14:22:50.678384"
$data.Range("F4").Value = "2021-10-05 14:22:50.678387"
$data.Range("F5").Value = "2021-10-05 14:22:50.678390"
$data.Range("F6").Value = "2021-10-05 14:22:50.678392"
$data.Range("F7").Value = "2021-10-05 14:22:50.678395"
$data.Range("F8").Value = "2021-10-05 14:22:50.678398"
$data.Range("F9").Value = "2021-10-05 14:22:50.678400"
$data.Range("F10").Value = "2021-10-05 14:22:50.678403"
$data.Range("F11").Value = "2021-10-05 14:22:50.678405"
$data.Range("F12").Value = "2021-10-05 14:22:50.678408"
$data.Range("F13").Value = "2021-10-05 14:22:50.678410"
$data.Range("F14").Value = "2021-10-05 14:22:50.678413"
$data.Range("F15").Value = "2021-10-05 14:22:50.678415"
$data.Range("F16").Value = "2021-10-05 14:22:50.678417"
$data.Range("F17").Value = "2021-10-05 14:22:50.678420"
$data.Range("F18").Value = "2021-10-05 14:22:50.678422"
$data.Range("F19").Value = "2021-10-05 14:22:50.678425"
$data.Range("F20").Value = "2021-10-05 14:22:50.678427"
$data.Range("F21").Value = "2021-10-05 14:22:50.678430"
$data.Range("F22").Value = "2021-10-05 14:22:50.678432"
$data.Range("F23").Value = "2021-10-05 14:22:50.678435"
$data.Range("F24").Value = "2021-10-05 14:22:50.678437"
$data.Range("F25").Value = "2021-10-05 14:22:50.678440"
$data.Range("F26").Value = "2021-10-05 14:22:50.678442"
$data.Range("F27").Value = "2021-10-05 14:22:50.678445"
$data.Range("F28").Value = "2021-10-05 14:22:50.678447"
$data.Range("F29").Value = "2021-10-05 14:22:50.678450"
$data.Range("F30").Value = "2021-10-05 14:22:50.678452"
$data.Range("F31").Value = "2021-10-05 14:22:50.678455"
$data.Range("F32").Value = "2021-10-05 14:22:50.678457"
$data.Range("F33").Value = "2021-10-05 14:22:50.678460"
$data.Range("F34").Value = "2021-10-05 14:22:50.678462"
$data.Range("F35").Value = "2021-10-05 14:22:50.678465"
$data.Range("F36").Value = "2021-10-05 14:22:50.678467"
$data.Range("F37").Value = "2021-10-05 14:22:50.678470"
$data.Range("F38").Value = "2021-10-05 14:22:50.678472"
$data.Range("F39").Value = "2021-10-05 14:22:50.678475"
$data.Range("F40").Value = "2021-10-05 14:22:50.678477"
$data.Range("F41").Value = "2021-10-05 14:22:50.678480"
$data.Range("F42").Value = "2021-10-05 14:22:50.678482"
$data.Range("F43").Value = "2021-10-05 14:22:50.678485"
$data.Range("F44").Value = "2021-10-05 14:22:50.678487"
$data.Range("F45").Value = "2021-10-05 14:22:50.678490"
$data.Range("F46").Value = "2021-10-05 14:22:50.678492"
$data.Range("F47").Value = "2021-10-05 14:22:50.678494"
$data.Range("F48").Value = "2021-10-05 14:22:50.678497"
$data.Range("F49").Value = "2021-10-05 14:22:50.678499"
$data.Range("F50").Value = "2021-10-05 14:22:50.678502"
$data.Range("F51").Value = "2021-10-05 14:22:50.678504"
$data.Range("F52").Value = "2021-10-05 14:22:50.678507"
$data.Range("F53").Value = "2021-10-05 14:22:50.678509"
$data.Range("F54").Value = "2021-10-05 14:22:50.678512"
$data.Range("F55").Value = "2021-10-05 14:22:50.678514"
$data.Range("F56").Value = "2021-10-05 14:22:50.678517"
$data.Range("F57").Value = "2021-10-05 14:22:50.678519"
$data.Range("F58").Value = "2021-10-05 14:22:50.678522"
$data.Range("F59").Value = "2021-10-05 14:22:50.678524"
$data.Range("F60").Value = "2021-10-05 14:22:50.678527"
$data.Range("F61").Value = "2021-10-05 14:22:50.678529"
$data.Range("F62").Value = "2021-10-05 14:22:50.678532"
$data.Range("F63").Value = "2021-10-05 14:22:50.678534"
$data.Range("F64").Value = "2021-10-05 14:22:50.678537"
$data.Range("F65").Value = "2021-10-05 14:22:50.678539"
$data.Range("F66").Value = "2021-10-05 14:22:50.678543"
$data.Range("F67").Value = "2021-10-05 14:22:50.678545"
$data.Range("F68").Value = "2021-10-05 14:22:50.678548"
$data.Range("F69").Value = "2021-10-05 14:22:50.678550"
$data.Range("F70").Value = "2021-10-05 14:22:50.678553"
$data.Range("F71").Value = "2021-10-05 14:22:50.678555"
$data.Range("F72").Value = "2021-10-05 14:22:50.678558"
$data.Range("F73").Value = "2021-10-05 14:22:50.678560"
$data.Range("F74").Value = "2021-10-05 14:22:50.678562"
$data.Range("F75").Value = "2021-10-05 14:22:50.678565"
$data.Range("F76").Value = "2021-10-05 14:22:50.678567"
$data.Range("F77").Value = "2021-10-05 14:22:50.678570"

# --- Add a dedicated "metadata" worksheet, placed after "data" -------------
$metadata = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$metadata.Name = "metadata"
# Match the outline defaults ("summary rows below detail" / "summary
# columns right of detail") used on the "data" sheet.
$metadata.Outline.SummaryRow = 1
$metadata.Outline.SummaryColumn = 1

# Header row — same bold/border/centered header style used on "data" row 1.
$data.Range("B1:F1").Copy()
$metadata.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"

$data.Range("F1").Copy()
$metadata.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$metadata.Range("G1").Value = "panel_get_request"

# Data row 2 — index column carries the same style as the data sheet's A2.
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$metadata.Range("A2").Value = 0

$metadata.Range("B2").Value = "Structural basal ganglia disorders"
$metadata.Range("C2").Value = 180
# Keep "1.19" as literal text (not the number 1.19) without tainting the
# cell's style: build it as a text formula, then flatten to a static value.
$metadata.Range("D2").Formula = '="1.19"'
$metadata.Range("D2").Copy()
$metadata.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$metadata.Range("E2").Value = "2021-08-20T13:49:07.993398Z"
$metadata.Range("F2").Value = "2021-10-05 14:22:50.674746"
$metadata.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/180/?format=json"

# Leave "data" as the active/selected sheet (unchanged from before the edit).
$data.Activate()
